$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly data row is inserted at row 8, pushing the existing
# rows 8-17 down to rows 9-18 (dimension grows from A1:T17 to A1:T18).
$ws.Rows.Item(8).Insert()

# Columns A-C, E-K are identical for every data row on this sheet, so
# copy them straight from the row above (which, after the insert, still
# holds the old row-8 values) into the freshly inserted row 8.
$ws.Range("A7:C7").Copy()
$ws.Range("A8").PasteSpecial()
$ws.Range("E7:K7").Copy()
$ws.Range("E8").PasteSpecial()

$ws.Cells.Item(8, 1).Value = 1
$ws.Cells.Item(8, 4).Value = 45028
$ws.Cells.Item(8, 12).Value = "Segunda"
$ws.Cells.Item(8, 13).Value = 200
$ws.Cells.Item(8, 14).Value = 21000
$ws.Cells.Item(8, 15).Value = 22000
$ws.Cells.Item(8, 16).Value = 21500
$ws.Cells.Item(8, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(8, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(8, 19).Value = 1075
$ws.Cells.Item(8, 20).Value = 20
